$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row at position 20 (pushes old row 20.. down by one;
#    this relocates the old signature-block rows 24/25 to 25/26).
# ------------------------------------------------------------------
$ws.Rows("20").Insert()

# ------------------------------------------------------------------
# 2) Duplicate row 19 (the former "last data row", which carries the
#    bottom-border look) down into the newly created row 20, formatting
#    and values both, as a base to edit from.
# ------------------------------------------------------------------
$ws.Range("B19:J19").Copy($ws.Range("B20:J20"))

# ------------------------------------------------------------------
# 3) Row 19 becomes a normal (non-last) data row, so give it the same
#    formatting as row 18 (the previous normal row immediately above).
# ------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) Write the final period / Valor Mora data, 2504 .. 2508, ascending.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2504"
$ws.Range("F16").Value = 45552

$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56940

$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 56940

$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940

$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940

# ------------------------------------------------------------------
# 5) Update the summary figures: total Valor Mora and period count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 273312
$ws.Range("F13").Value = 5
